$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 2758
$ws.Range("F9").Value = 6273
$ws.Range("F11").Value = 75
$ws.Range("F13").Value = 4996
$ws.Range("F15").Value = 542
$ws.Range("F16").Value = 2607
$ws.Range("F17").Value = 1343
$ws.Range("F18").Value = 1513
$ws.Range("F20").Value = 305
$ws.Range("F21").Value = 119
$ws.Range("F23").Value = 1056
$ws.Range("F24").Value = 234
$ws.Range("F26").Value = 526
$ws.Range("F27").Value = 1365
$ws.Range("F28").Value = 1029
$ws.Range("F29").Value = 2092
$ws.Range("F30").Value = 307
$ws.Range("F32").Value = 19
$ws.Range("F33").Value = 21
$ws.Range("F34").Value = 84
$ws.Range("F35").Value = 246
$ws.Range("F36").Value = 1491
$ws.Range("F42").Value = 290
$ws.Range("F43").Value = 2260
$ws.Range("F44").Value = 2540
$ws.Range("F46").Value = 128
$ws.Range("F47").Value = 270
$ws.Range("F49").Value = 91

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 404
$ws.Range("F8").Value = 320
$ws.Range("F10").Value = 88
$ws.Range("F11").Value = 202
$ws.Range("F12").Value = 6
$ws.Range("F23").Value = 367
$ws.Range("F31").Value = 5
$ws.Range("F36").Value = 25
$ws.Range("F37").Value = 26

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 1695
$ws.Range("F7").Value = 567
$ws.Range("F8").Value = 1479
$ws.Range("F9").Value = 1810
$ws.Range("F10").Value = 2500
$ws.Range("F11").Value = 836
$ws.Range("F12").Value = 716
$ws.Range("F13").Value = 4

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 567
$ws.Range("F7").Value = 2758
$ws.Range("F9").Value = 1479
$ws.Range("F11").Value = 2500
$ws.Range("F12").Value = 6273
$ws.Range("F13").Value = 836
$ws.Range("F16").Value = 4996
$ws.Range("F17").Value = 2607
$ws.Range("F18").Value = 1343
$ws.Range("F19").Value = 1513
$ws.Range("F21").Value = 119
$ws.Range("F23").Value = 320
$ws.Range("F24").Value = 234
$ws.Range("F25").Value = 88
$ws.Range("F27").Value = 1365
$ws.Range("F28").Value = 1029
$ws.Range("F29").Value = 2092
$ws.Range("F30").Value = 307
$ws.Range("F32").Value = 19
$ws.Range("F34").Value = 21
$ws.Range("F35").Value = 246
$ws.Range("F42").Value = 290
$ws.Range("F44").Value = 2260
$ws.Range("F45").Value = 2540
$ws.Range("F46").Value = 128
$ws.Range("F47").Value = 270
